$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F2").Value = 3.4
$ws.Range("G2").Value = 4.6
$ws.Range("H2").Value = 2.02
$ws.Range("I2").Value = 2.44
$ws.Range("J2").Value = 3.3
$ws.Range("M2").Value = 1.02
$ws.Range("N2").Value = 1.76
$ws.Range("P2").Value = 1.76
$ws.Range("Q2").Value = 1.9
$ws.Range("R2").Value = 1.25
$ws.Range("S2").Value = 1.9
$ws.Range("V2").Value = 1.73
$ws.Range("W2").Value = 1.28
$ws.Range("G3").Value = 1.69
$ws.Range("H3").Value = 2.44
$ws.Range("P3").Value = 1.81
$ws.Range("Q3").Value = 1.71
$ws.Range("F4").Value = 1.97
$ws.Range("G4").Value = 1.99
$ws.Range("I4").Value = 4.4
$ws.Range("K4").Value = 3.9
$ws.Range("Q4").Value = 1.87
$ws.Range("T4").Value = 1.79
$ws.Range("AE4").Value = 180
$ws.Range("H5").Value = 12.5
$ws.Range("J5").Value = 5.3
$ws.Range("O5").Value = 1.31
$ws.Range("P5").Value = 2
$ws.Range("Q5").Value = 1.95
$ws.Range("R5").Value = 1.37
$ws.Range("X5").Value = 17
$ws.Range("Y5").Value = 34
$ws.Range("Z5").Value = 140
$ws.Range("AB5").Value = 7
$ws.Range("AC5").Value = 12.5
$ws.Range("AD5").Value = 55
$ws.Range("AE5").Value = 370
$ws.Range("AF5").Value = 7
$ws.Range("AG5").Value = 11
$ws.Range("AH5").Value = 40
$ws.Range("AJ5").Value = 10
$ws.Range("AK5").Value = 17
$ws.Range("AL5").Value = 60
$ws.Range("AM5").Value = 350
$ws.Range("AN5").Value = 6.6
$ws.Range("F6").Value = 2.6
$ws.Range("G6").Value = 3.5
$ws.Range("H6").Value = 2.42
$ws.Range("I6").Value = 3.1
$ws.Range("J6").Value = 3.15
$ws.Range("Q6").Value = 1.71
$ws.Range("F7").Value = 1.55
$ws.Range("G7").Value = 1.96
$ws.Range("H7").Value = 3.75
$ws.Range("I7").Value = 7
$ws.Range("J7").Value = 3.85
$ws.Range("K7").Value = 11
$ws.Range("N7").Value = 1.93
$ws.Range("P7").Value = 1.92
$ws.Range("R7").Value = 1.38
$ws.Range("S7").Value = 2.52
$ws.Range("F8").Value = 10.5
$ws.Range("G8").Value = 11.5
$ws.Range("H8").Value = 1.29
$ws.Range("I8").Value = 1.3
$ws.Range("Q8").Value = 1.29
$ws.Range("T8").Value = 1.62
$ws.Range("U8").Value = 2.54
$ws.Range("X8").Value = 60
$ws.Range("AC8").Value = 18.5
$ws.Range("AD8").Value = 12
$ws.Range("F9").Value = 1.9
$ws.Range("G9").Value = 1.93
$ws.Range("H9").Value = 4.2
$ws.Range("P9").Value = 2.44
$ws.Range("R9").Value = 1.57
$ws.Range("S9").Value = 2.64
$ws.Range("T9").Value = 1.65
$ws.Range("U9").Value = 2.44
$ws.Range("AG9").Value = 10.5
$ws.Range("AN9").Value = 9
$ws.Range("F10").Value = 1.84
$ws.Range("G10").Value = 1.85
$ws.Range("H10").Value = 4.4
$ws.Range("I10").Value = 4.5
$ws.Range("J10").Value = 4.2
$ws.Range("X10").Value = 24
$ws.Range("Y10").Value = 22
$ws.Range("Z10").Value = 38
$ws.Range("AD10").Value = 19
$ws.Range("AH10").Value = 17
$ws.Range("AK10").Value = 18.5
$ws.Range("AM10").Value = 75
$ws.Range("AN10").Value = 9
$ws.Range("AO10").Value = 42
$ws.Range("J11").Value = 9.4
$ws.Range("K11").Value = 10
$ws.Range("G14").Value = 1.55
$ws.Range("H14").Value = 6.4
$ws.Range("I14").Value = 15.5
$ws.Range("K14").Value = 11
$ws.Range("N14").Value = 2.14
$ws.Range("O14").Value = 1.19
$ws.Range("P14").Value = 2.14
$ws.Range("Q14").Value = 1.56
$ws.Range("R14").Value = 1.46
$ws.Range("S14").Value = 2.3
$ws.Range("W14").Value = 2.82
